$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 68.664
$ws.Range("D2").Value = 68.664
$ws.Range("E2").Value = 2.47645532
$ws.Range("F2").Value = 0.00086356
$ws.Range("G2").Value = 0.05921297
$ws.Range("H2").Value = 4.08556316
$ws.Range("I2").Value = 5.334375273596671
$ws.Range("J2").Value = 5.334375273596671
$ws.Range("K2").Value = 0.1968626769498817
$ws.Range("L2").Value = 0.00006962883436173964
$ws.Range("M2").Value = 0.005940344347768912
$ws.Range("N2").Value = 0.649636654064962
$ws.Range("C3").Value = 88.009
$ws.Range("D3").Value = 88.009
$ws.Range("E3").Value = 1.95433027
$ws.Range("F3").Value = 0.00059699
$ws.Range("G3").Value = 0.05203552000000001
$ws.Range("H3").Value = 4.616124459999999
$ws.Range("I3").Value = 11.76236678021359
$ws.Range("J3").Value = 11.76236678021359
$ws.Range("K3").Value = 0.2595545276207594
$ws.Range("L3").Value = 0.00006070358164741868
$ws.Range("M3").Value = 0.004749481667056491
$ws.Range("N3").Value = 0.9561897840302375
$ws.Range("C4").Value = 38.825
$ws.Range("D4").Value = 77.559
$ws.Range("E4").Value = 2.21976254
$ws.Range("F4").Value = 0.00162311
$ws.Range("G4").Value = 0.03139647
$ws.Range("H4").Value = 1.23813492
$ws.Range("I4").Value = 5.465734099757255
$ws.Range("J4").Value = 10.89914264106633
$ws.Range("K4").Value = 0.2956742843521831
$ws.Range("L4").Value = 0.0001082855191318472
$ws.Range("M4").Value = 0.004017594741459945
$ws.Range("N4").Value = 0.3309608445411829
$ws.Range("C5").Value = 48.977
$ws.Range("D5").Value = 95.48099999999999
$ws.Range("E5").Value = 1.80847664
$ws.Range("F5").Value = 0.00113542
$ws.Range("G5").Value = 0.02754358000000001
$ws.Range("H5").Value = 1.36982133
$ws.Range("I5").Value = 7.730474191874081
$ws.Range("J5").Value = 14.13107661596689
$ws.Range("K5").Value = 0.26521108367037
$ws.Range("L5").Value = 0.0001234109146463709
$ws.Range("M5").Value = 0.003707729256545013
$ws.Range("N5").Value = 0.3807253529087321
$ws.Range("C6").Value = 24.016
$ws.Range("D6").Value = 95.84099999999999
$ws.Range("E6").Value = 1.80853719
$ws.Range("F6").Value = 0.00309052
$ws.Range("G6").Value = 0.01856005
$ws.Range("H6").Value = 0.45732042
$ws.Range("I6").Value = 3.864541198988397
$ws.Range("J6").Value = 15.42457897761351
$ws.Range("K6").Value = 0.2906751818496847
$ws.Range("L6").Value = 0.0002616596938685731
$ws.Range("M6").Value = 0.003398262185791936
$ws.Range("N6").Value = 0.155548288691771
$ws.Range("C7").Value = 28.025
$ws.Range("D7").Value = 102.145
$ws.Range("E7").Value = 1.69266166
$ws.Range("F7").Value = 0.00215418
$ws.Range("G7").Value = 0.01493945
$ws.Range("H7").Value = 0.42923915
$ws.Range("I7").Value = 5.411807707472861
$ws.Range("J7").Value = 15.56445396495577
$ws.Range("K7").Value = 0.2564875004446195
$ws.Range("L7").Value = 0.0002280533620360286
$ws.Range("M7").Value = 0.002423054516656753
$ws.Range("N7").Value = 0.1504849184478556
$ws.Range("C8").Value = 16.766
$ws.Range("D8").Value = 100.331
$ws.Range("E8").Value = 1.77649457
$ws.Range("F8").Value = 0.00465617
$ws.Range("G8").Value = 0.01299438
$ws.Range("H8").Value = 0.2278982
$ws.Range("I8").Value = 3.636546228012023
$ws.Range("J8").Value = 21.76463043881375
$ws.Range("K8").Value = 0.4549831109261563
$ws.Range("L8").Value = 0.0003660237758732409
$ws.Range("M8").Value = 0.002942759223492035
$ws.Range("N8").Value = 0.09819261343753762
$ws.Range("C9").Value = 19.693
$ws.Range("D9").Value = 100.229
$ws.Range("E9").Value = 1.73325661
$ws.Range("F9").Value = 0.00309312
$ws.Range("G9").Value = 0.01006144
$ws.Range("H9").Value = 0.20524318
$ws.Range("I9").Value = 4.253568365953913
$ws.Range("J9").Value = 16.65124087937061
$ws.Range("K9").Value = 0.2919522439916838
$ws.Range("L9").Value = 0.0003542334639005001
$ws.Range("M9").Value = 0.002020930747820788
$ws.Range("N9").Value = 0.08598169339159405
$ws.Range("C10").Value = 12.032
$ws.Range("D10").Value = 95.96299999999999
$ws.Range("E10").Value = 1.89867041
$ws.Range("F10").Value = 0.00600858
$ws.Range("G10").Value = 0.00902055
$ws.Range("H10").Value = 0.11569048
$ws.Range("I10").Value = 3.139240486300919
$ws.Range("J10").Value = 25.01782409648214
$ws.Range("K10").Value = 0.5621345139248056
$ws.Range("L10").Value = 0.0004348124255876843
$ws.Range("M10").Value = 0.002377026430680179
$ws.Range("N10").Value = 0.0587414584353224
$ws.Range("C11").Value = 15.069
$ws.Range("D11").Value = 93.55800000000001
$ws.Range("E11").Value = 1.87360093
$ws.Range("F11").Value = 0.00382647
$ws.Range("G11").Value = 0.00712754
$ws.Range("H11").Value = 0.11255008
$ws.Range("I11").Value = 3.755308421600194
$ws.Range("J11").Value = 18.02835704266081
$ws.Range("K11").Value = 0.3601807146122814
$ws.Range("L11").Value = 0.0004997885830305928
$ws.Range("M11").Value = 0.001663770140217536
$ws.Range("N11").Value = 0.05429500054404345
$ws.Range("C12").Value = 9.045999999999999
$ws.Range("D12").Value = 90.143
$ws.Range("E12").Value = 2.07618973
$ws.Range("F12").Value = 0.007664120000000001
$ws.Range("G12").Value = 0.0069627
$ws.Range("H12").Value = 0.06916602000000001
$ws.Range("I12").Value = 2.800665586557547
$ws.Range("J12").Value = 27.89427614478715
$ws.Range("K12").Value = 0.7030763380091064
$ws.Range("L12").Value = 0.0005826775598379765
$ws.Range("M12").Value = 0.002272251710932976
$ws.Range("N12").Value = 0.04302386307333819
$ws.Range("C13").Value = 12.052
$ws.Range("D13").Value = 84.774
$ws.Range("E13").Value = 2.07545842
$ws.Range("F13").Value = 0.00450943
$ws.Range("G13").Value = 0.00537033
$ws.Range("H13").Value = 0.06835468
$ws.Range("I13").Value = 3.180788174873236
$ws.Range("J13").Value = 17.14581298460073
$ws.Range("K13").Value = 0.41890826142342
$ws.Range("L13").Value = 0.0006953722151607854
$ws.Range("M13").Value = 0.002423054516656753
$ws.Range("N13").Value = 0.03632059136984219
